# Apply updates described by the commit:
# "Updated detailed indicator quantile results."
#
# 1) Update the shared "ScriptLatestRunVersion" text (Git Commit ID) used in
#    column AJ (rows 2-80) from the old commit hash to the new one.
# 2) Update the "pid" values in column AH (rows 2-80) from 25080 to 21528.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$oldCommitText = "IndicatorQuantiles.R, Git Commit ID: 2e3ff9a54734c37c56b32bb788c6f054c2509b6b"
$newCommitText = "IndicatorQuantiles.R, Git Commit ID: db49f0f869e1f5a8558dc746458075a467cf2c41"

# Find the last used row in column A to know how far data extends.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 80 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ajCell = $ws.Cells.Item($r, 36)  # column AJ = 36
    if ($ajCell.Value2 -eq $oldCommitText) {
        $ajCell.Value2 = $newCommitText
    }

    $ahCell = $ws.Cells.Item($r, 34)  # column AH = 34
    if ($ahCell.Value2 -eq 25080) {
        $ahCell.Value2 = 21528
    }
}
